{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of 1-based row index (document order) -> new text for that row's\n// single cell. These correspond 1:1 to the <w:tr> rows in the single\n// column results table.\nconst rowUpdates = {\n  1: \"0M\",\n  2: \"0M\",\n  3: \"0M\",\n  4: \"44\",\n  6: \"0.00042\",\n  7: \"0.00017\",\n  8: \"0.00005\",\n  9: \"0.00021\",\n  10: \"0.00023\",\n  11: \"0.00034\",\n  12: \"0.00761\",\n  44: \"99.94\",\n  45: \"0.01\",\n  46: \"12\",\n};\n\n// First load the first paragraph of every target cell's body so we can\n// replace its text range while preserving the run formatting.\nconst paragraphsByRow = {};\nfor (const rowIdx of Object.keys(rowUpdates)) {\n  const cell = table.getCell(parseInt(rowIdx, 10) - 1, 0);\n  const body = cell.body;\n  body.paragraphs.load(\"items\");\n  paragraphsByRow[rowIdx] = body;\n}\nawait context.sync();\n\nfor (const [rowIdx, newText] of Object.entries(rowUpdates)) {\n  const body = paragraphsByRow[rowIdx];\n  const para = body.paragraphs.items[0];\n  const range = para.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Map of 1-based row index (document order) -> new text for that row's\n# single cell. These correspond 1:1 to the <w:tr> rows in the single\n# column results table.\n$rowUpdates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"44\"\n    6  = \"0.00042\"\n    7  = \"0.00017\"\n    8  = \"0.00005\"\n    9  = \"0.00021\"\n    10 = \"0.00023\"\n    11 = \"0.00034\"\n    12 = \"0.00761\"\n    44 = \"99.94\"\n    45 = \"0.01\"\n    46 = \"12\"\n}\n\nforeach ($rowIdx in $rowUpdates.Keys) {\n    $t.Cell($rowIdx, 1).Range.Text = $rowUpdates[$rowIdx]\n}\n"}
